$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new trade row (row 12) with the same column layout as the
# existing rows (A:H).
$ws.Range("A12").Value = 9219.31
$ws.Range("B12").Value = 9138.89
$ws.Range("C12").Value = 105.78
$ws.Range("D12").Value = 106.71
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 0.88

# G column holds a date/time value (formatted like the rest of column G,
# style index 1 / numFmtId 22). Copy the existing date format from the row
# above so the cell reuses the same style instead of minting a new one.
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").Value = 42620.766215277778

$ws.Range("H12").Value = $true
